$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows from bottom to top so row indices of the ones still to be
# removed remain valid.
$ws.Rows.Item(28).EntireRow.Delete()
$ws.Rows.Item(26).EntireRow.Delete()
